# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers (AD1:AF1) -----------------------------------------------------
# Copy the formatting of the last existing header cell (AC1, style s="1":
# bold font, thin border, centered/top-aligned) onto the three new header
# cells before writing their text, so they match the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Season-record data (rows 2-48) ---------------------------------------
$lastRow = 48
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 66   # column AD - Wins
    $ws.Cells.Item($row, 31).Value = 96   # column AE - Losses
    $ws.Cells.Item($row, 32).Value = 0    # column AF - Ties
}
